$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.1424201652695178
$ws.Cells.Item(2, 8).Value = 50.35772395746098
$ws.Cells.Item(2, 9).Value = 1.165065652347801
$ws.Cells.Item(3, 7).Value = 0.1519146242324391
$ws.Cells.Item(3, 8).Value = 129.8466419588245
$ws.Cells.Item(4, 7).Value = 0.01626309980923994
$ws.Cells.Item(4, 8).Value = 73.46811390566157
$ws.Cells.Item(5, 7).Value = 0.02379283006678161
$ws.Cells.Item(5, 8).Value = 277.693528930069
$ws.Cells.Item(6, 7).Value = -0.2272541810922561
$ws.Cells.Item(6, 8).Value = -2.73897431696032
$ws.Cells.Item(7, 7).Value = -0.2587639806376527
$ws.Cells.Item(7, 8).Value = -3.557083191018312
$ws.Cells.Item(8, 7).Value = -0.3743793175321739
$ws.Cells.Item(8, 8).Value = -1.11367519582157
$ws.Cells.Item(9, 7).Value = -0.462289691914912
$ws.Cells.Item(9, 8).Value = -15.96834771074314
$ws.Cells.Item(10, 7).Value = 0.004862212323446856
$ws.Cells.Item(10, 8).Value = -69.99539741803258
$ws.Cells.Item(11, 7).Value = 0.0179918748571623
$ws.Cells.Item(11, 8).Value = 211.9881951442589
$ws.Cells.Item(12, 7).Value = 0.2376503796500293
$ws.Cells.Item(12, 8).Value = 4.604795012735252
$ws.Cells.Item(13, 7).Value = 0.2560882449064057
$ws.Cells.Item(13, 8).Value = -2.753977352883232
$ws.Cells.Item(14, 7).Value = -0.01639610553023982
$ws.Cells.Item(14, 8).Value = -71.37701175776884
$ws.Cells.Item(15, 7).Value = 0.01496293281258154
$ws.Cells.Item(15, 8).Value = -25.87733381269883
$ws.Cells.Item(16, 7).Value = 0.1421679474362366
$ws.Cells.Item(16, 8).Value = 20.48222483856073
$ws.Cells.Item(17, 7).Value = 0.2057615078853936
$ws.Cells.Item(17, 8).Value = -5.973239490440858
$ws.Cells.Item(18, 7).Value = 0.0475624122425268
$ws.Cells.Item(18, 8).Value = -21.34102214183326
$ws.Cells.Item(19, 7).Value = 0.06977307085352141
$ws.Cells.Item(19, 8).Value = -22.55121565822473
$ws.Cells.Item(20, 7).Value = -0.1549612290636011
$ws.Cells.Item(20, 8).Value = -6.479894833216439
$ws.Cells.Item(21, 7).Value = -0.1967736510113334
$ws.Cells.Item(21, 8).Value = 1.529422208367743
$ws.Cells.Item(22, 7).Value = 0.04881641967134771
$ws.Cells.Item(22, 8).Value = -10.24453566373104
$ws.Cells.Item(23, 7).Value = 0.02561024320794179
$ws.Cells.Item(23, 8).Value = -37.29149813838335
$ws.Cells.Item(24, 7).Value = 0.1302089078626761
$ws.Cells.Item(24, 8).Value = 12.5067347632109
$ws.Cells.Item(25, 7).Value = 0.1720563124983373
$ws.Cells.Item(25, 8).Value = 13.13835198880915
$ws.Cells.Item(26, 7).Value = 0.02882637362060067
$ws.Cells.Item(26, 8).Value = -45.48387794770135
$ws.Cells.Item(27, 7).Value = 0.01290957487798713
$ws.Cells.Item(27, 8).Value = -74.4196720142334
$ws.Cells.Item(28, 7).Value = 0.1778414196775265
$ws.Cells.Item(28, 8).Value = 16.30382078824532
$ws.Cells.Item(29, 7).Value = 0.1939093667549562
$ws.Cells.Item(29, 8).Value = 13.59452784430934
$ws.Cells.Item(30, 7).Value = -0.00303138070973302
$ws.Cells.Item(30, 8).Value = -115.4930351121823
$ws.Cells.Item(31, 7).Value = -0.006079593786075354
$ws.Cells.Item(31, 8).Value = -162.6441236688637
$ws.Cells.Item(32, 7).Value = 0.02174531191890518
$ws.Cells.Item(32, 8).Value = -41.6910595962214
$ws.Cells.Item(33, 7).Value = 0.02454515545609308
$ws.Cells.Item(33, 8).Value = -5.970368433940076
$ws.Cells.Item(34, 7).Value = 0.1036443620529278
$ws.Cells.Item(34, 8).Value = -19.00759950607765
$ws.Cells.Item(35, 7).Value = 0.1191880133549298
$ws.Cells.Item(35, 8).Value = -7.36245611223776
$ws.Cells.Item(36, 7).Value = -0.04245097173126619
$ws.Cells.Item(36, 8).Value = -382.4061305241412
$ws.Cells.Item(37, 7).Value = 0.007847571375231035
$ws.Cells.Item(37, 8).Value = -48.75720841059535
$ws.Cells.Item(38, 7).Value = -0.04048237637119766
$ws.Cells.Item(38, 8).Value = -1882.085146885811
$ws.Cells.Item(39, 7).Value = 0.00551830619617054
$ws.Cells.Item(39, 8).Value = 116.5173190271454
$ws.Cells.Item(40, 7).Value = 0.1550405278804526
$ws.Cells.Item(40, 8).Value = 5.077320898199825
$ws.Cells.Item(41, 7).Value = 0.1345867270597043
$ws.Cells.Item(41, 8).Value = -16.61252226345421
$ws.Cells.Item(42, 7).Value = 0.06084472567608996
$ws.Cells.Item(42, 8).Value = -5.76192060942398
$ws.Cells.Item(43, 7).Value = 0.05019213739977633
$ws.Cells.Item(43, 8).Value = 44.39407922810257
$ws.Cells.Item(44, 7).Value = 0.0175896988863587
$ws.Cells.Item(44, 8).Value = 24.63664536274487
$ws.Cells.Item(45, 7).Value = 0.01517207453426739
$ws.Cells.Item(45, 8).Value = -63.04724946818273
$ws.Cells.Item(46, 7).Value = -0.04730867865013868
$ws.Cells.Item(46, 8).Value = 28.1240103427208
$ws.Cells.Item(47, 7).Value = -0.04220948004513633
$ws.Cells.Item(47, 8).Value = -2.178125654373262
$ws.Cells.Item(48, 7).Value = -0.109997798354968
$ws.Cells.Item(48, 8).Value = 12.68311898417387
$ws.Cells.Item(49, 7).Value = -0.1240440186921968
$ws.Cells.Item(49, 8).Value = 37.18698846810843
$ws.Cells.Item(50, 7).Value = 0.1053878828021463
$ws.Cells.Item(50, 8).Value = -3.201660864794182
$ws.Cells.Item(51, 7).Value = 0.08093681869259406
$ws.Cells.Item(51, 8).Value = -19.28178909175136
$ws.Cells.Item(52, 7).Value = 0.0502435763729707
$ws.Cells.Item(52, 8).Value = -15.7256386705489
$ws.Cells.Item(53, 7).Value = 0.0503854953393062
$ws.Cells.Item(53, 8).Value = -25.40199025589368
$ws.Cells.Item(54, 7).Value = -0.09246111151462176
$ws.Cells.Item(54, 8).Value = -32.23733424189755
$ws.Cells.Item(55, 7).Value = -0.08397365933765766
$ws.Cells.Item(55, 8).Value = -8.723184905940586
$ws.Cells.Item(56, 7).Value = 0.09657761222459652
$ws.Cells.Item(56, 8).Value = 110.7481716138621
$ws.Cells.Item(57, 7).Value = 0.1155617708395022
$ws.Cells.Item(57, 8).Value = 2135.117242185928
